# Player-valuation-model: add a "Tuned Lasso Regression" results column (H)
# to the Model Performance Metrics sheet, matching the formats already used
# by the neighbouring "OLS"/"Lasso Regression - Untuned" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column header (H2) + metric values (H3:H8)
# ---------------------------------------------------------------------
$ws.Range("H2").Value = "Tuned Lasso Regression"

$ws.Range("H3").Value = 0.4277
$ws.Range("H4").Value = 0.4455
$ws.Range("H5").Value = 0.6674
$ws.Range("H6").Value = 1.6248
$ws.Range("H7").Value = 0.296
$ws.Range("H8").Value = 0.6151

# ---------------------------------------------------------------------
# 2. Formatting for the new column - copy from the analogous cells used
#    for the other model columns (D/F header style, F data-row style,
#    G percent style for the R-Squared row).
# ---------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("F4").Copy()
$ws.Range("H4").PasteSpecial(-4122)

$ws.Range("F5").Copy()
$ws.Range("H5").PasteSpecial(-4122)

$ws.Range("F6").Copy()
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("H7").PasteSpecial(-4122)

$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Column widths - columns resized (existing B-E shrank a bit) and the
#    new column H sized/bestfit like the rest of the metric columns.
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 23.333333333333336
$ws.Columns("C").ColumnWidth = 29.5
$ws.Columns("D").ColumnWidth = 25.5
$ws.Columns("E").ColumnWidth = 33.5
$ws.Columns("H").ColumnWidth = 25.166666666666664

# ---------------------------------------------------------------------
# 4. Sheet view / selection - user scrolled back to A1 and left the
#    selection on I15 instead of the old F1:G1 / topLeftCell=E1 state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("I15").Select()
